$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hospital")

# Update hospital name and team name (Kampong Cham -> Siem Reap)
$ws.Range("B2").Value = "Siem Reap Provincial Referral Hospital"
$ws.Range("B5").Value = "Siem Reap Microbiology Team"

# Update the active selection on the "hospital" sheet
$ws.Activate()
$ws.Range("C17").Select()
